$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "Emberlee"
$ws.Range("D16").Value = "Create Mocap Poses and put them in unity"

$ws.Range("D17").Select()
